$wb = $excel.ActiveWorkbook

# --- Update "Zadatak 1" sheet with corrected MPI program times ---
$ws1 = $wb.Worksheets.Item("Zadatak 1")

$ws1.Range("A3").Value = 0.000993
$ws1.Range("C3").Value = 0.008646
$ws1.Range("E3").Value = 0.177077
$ws1.Range("G3").Value = 4.21437

$ws1.Range("A4").Value = 0.00286
$ws1.Range("C4").Value = 0.008366
$ws1.Range("E4").Value = 0.176766
$ws1.Range("G4").Value = 4.212673

# Make "Zadatak 1" the active sheet/tab with the selection on G4
$ws1.Activate()
$ws1.Range("G4").Select()
